$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("DATA")
$wsDict = $wb.Worksheets.Item("DICTIONARY")

# --- DATA sheet (sheet1) changes --------------------------------------
# The model now executes at a segment level instead of account level, so
# the number of simulations (column K / "N") increases from 10 to 100 for
# every segment row, and the "Z-Non-Cyclical" segment's x0/dx0 values are
# updated. The x0 cell (E3) no longer derives its value from dx0 (G3) via
# a formula - it is now an independent static value.
$wsData.Range("K2").Value = 100
$wsData.Range("K3").Value = 100
$wsData.Range("K4").Value = 100
$wsData.Range("K5").Value = 100

$wsData.Range("G3").Value = 0.025
$wsData.Range("E3").Value = 0.025

# --- Sheet view / active sheet changes --------------------------------
# DATA becomes the active / selected sheet (it used to be DICTIONARY),
# with the selected cell now E3 instead of E12.
$wsData.Range("E3").Select()
$wsData.Activate()

# DICTIONARY keeps its own selection (H11) but is no longer the active tab.
$wsDict.Range("H11").Select()

# DICTIONARY picks up an explicit page setup (A4 / portrait).
$wsDict.PageSetup.PaperSize = 9
$wsDict.PageSetup.Orientation = 1

# Re-activate DATA last so it is the tab shown/selected when the workbook
# is reopened.
$wsData.Activate()
